$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3042.5757
$ws.Range("I64").Value = 2962.3333
$ws.Range("K64").Value = 2962.3333
$ws.Range("M64").Value = -2714.3333
$ws.Range("H67").Value = 3042.5757
$ws.Range("I67").Value = 2962.3333
$ws.Range("K67").Value = 2962.3333
$ws.Range("M67").Value = -2104.3333
$ws.Range("H76").Value = 3099.7727
$ws.Range("I76").Value = 3060.7334
$ws.Range("J76").Value = 3183.4285
$ws.Range("K76").Value = 3060.7334
$ws.Range("L76").Value = 3183.4285
$ws.Range("M76").Value = -2745.7334
$ws.Range("N76").Value = -3813.4285
$ws.Range("H79").Value = 3099.7727
$ws.Range("I79").Value = 3060.7334
$ws.Range("J79").Value = 3183.4285
$ws.Range("K79").Value = 3060.7334
$ws.Range("L79").Value = 3183.4285
$ws.Range("M79").Value = -1968.7334
$ws.Range("N79").Value = -5367.4285
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H138").Value = 2455.93
$ws.Range("J138").Value = 3560.1091
$ws.Range("L138").Value = 10680.3273
$ws.Range("N138").Value = -20960.3273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 2341.3333
$ws.Range("I39").Value = 809.6
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 809.6
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -289.6
$ws.Range("N39").Value = -11040
$ws.Range("H61").Value = 2314.5908
$ws.Range("I61").Value = 2018.4166
$ws.Range("K61").Value = 2018.4166
$ws.Range("M61").Value = -1806.4166
$ws.Range("H63").Value = 2568.3333
$ws.Range("I63").Value = 2568.3333
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2568.3333
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1882.3333
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2568.3333
$ws.Range("I66").Value = 2568.3333
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12841.6665
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9409.666499999999
$ws.Range("N66").ClearContents()
$ws.Range("H136").Value = 2314.5908
$ws.Range("I136").Value = 2018.4166
$ws.Range("K136").Value = 6055.2498
$ws.Range("M136").Value = -3505.2498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5716475.5
$ws.Range("I86").Value = 6668767.5
$ws.Range("J86").Value = 2722.8
$ws.Range("K86").Value = 6668767.5
$ws.Range("L86").Value = 2722.8
$ws.Range("M86").Value = -6667644.5
$ws.Range("N86").Value = -4968.8
$ws.Range("H89").Value = 5716475.5
$ws.Range("I89").Value = 6668767.5
$ws.Range("J89").Value = 2722.8
$ws.Range("K89").Value = 33343837.5
$ws.Range("L89").Value = 13614
$ws.Range("M89").Value = -33338221.5
$ws.Range("N89").Value = -24846
$ws.Range("H105").Value = 2185.6428
$ws.Range("I105").Value = 1708.7778
$ws.Range("J105").Value = 3044
$ws.Range("K105").Value = 1708.7778
$ws.Range("L105").Value = 3044
$ws.Range("M105").Value = 38.22219999999993
$ws.Range("N105").Value = -6538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2740.4324
$ws.Range("I31").Value = 1470.75
$ws.Range("J31").Value = 3707.8096
$ws.Range("K31").Value = 1470.75
$ws.Range("L31").Value = 3707.8096
$ws.Range("M31").Value = -1175.75
$ws.Range("N31").Value = -4297.809600000001
$ws.Range("H34").Value = 2740.4324
$ws.Range("I34").Value = 1470.75
$ws.Range("J34").Value = 3707.8096
$ws.Range("K34").Value = 1470.75
$ws.Range("L34").Value = 3707.8096
$ws.Range("M34").Value = -1268.75
$ws.Range("N34").Value = -4111.809600000001
$ws.Range("H58").Value = 1671.836
$ws.Range("I58").Value = 649.0278
$ws.Range("J58").Value = 3144.68
$ws.Range("K58").Value = 649.0278
$ws.Range("L58").Value = 3144.68
$ws.Range("M58").Value = -446.0278
$ws.Range("N58").Value = -3550.68
$ws.Range("H62").Value = 2933.1428
$ws.Range("I62").Value = 2438.75
$ws.Range("J62").Value = 3237.3845
$ws.Range("K62").Value = 2438.75
$ws.Range("L62").Value = 3237.3845
$ws.Range("M62").Value = -1814.75
$ws.Range("N62").Value = -4485.3845
$ws.Range("H65").Value = 2933.1428
$ws.Range("I65").Value = 2438.75
$ws.Range("J65").Value = 3237.3845
$ws.Range("K65").Value = 12193.75
$ws.Range("L65").Value = 16186.9225
$ws.Range("M65").Value = -9073.75
$ws.Range("N65").Value = -22426.9225
$ws.Range("H134").Value = 2373.311
$ws.Range("I134").Value = 2720.8057
$ws.Range("J134").Value = 983.3333
$ws.Range("K134").Value = 8162.4171
$ws.Range("L134").Value = 2949.9999
$ws.Range("M134").Value = -5627.4171
$ws.Range("N134").Value = -8019.9999
$ws.Range("H136").Value = 1671.836
$ws.Range("I136").Value = 649.0278
$ws.Range("J136").Value = 3144.68
$ws.Range("K136").Value = 1947.0834
$ws.Range("L136").Value = 9434.039999999999
$ws.Range("M136").Value = 602.9166
$ws.Range("N136").Value = -14534.04

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1001.1539
$ws.Range("J5").Value = 1133.7894
$ws.Range("L5").Value = 3401.3682
$ws.Range("N5").Value = -3625.3682
$ws.Range("H120").Value = 13133
$ws.Range("I120").Value = 12999
$ws.Range("J120").Value = 13200
$ws.Range("K120").Value = 38997
$ws.Range("L120").Value = 39600
$ws.Range("M120").Value = -34159
$ws.Range("N120").Value = -49276
$ws.Range("H135").Value = 1001.1539
$ws.Range("J135").Value = 1133.7894
$ws.Range("L135").Value = 10204.1046
$ws.Range("N135").Value = -15274.1046

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 92160040
$ws.Range("I70").Value = 207354690
$ws.Range("J70").Value = 4320
$ws.Range("K70").Value = 207354690
$ws.Range("L70").Value = 4320
$ws.Range("M70").Value = -207354420
$ws.Range("N70").Value = -4860
$ws.Range("H73").Value = 92160040
$ws.Range("I73").Value = 207354690
$ws.Range("J73").Value = 4320
$ws.Range("K73").Value = 207354690
$ws.Range("L73").Value = 4320
$ws.Range("M73").Value = -207353754
$ws.Range("N73").Value = -6192
$ws.Range("H132").Value = 1776.661
$ws.Range("I132").Value = 1297.561
$ws.Range("J132").Value = 2867.9443
$ws.Range("K132").Value = 3892.683
$ws.Range("L132").Value = 8603.832900000001
$ws.Range("M132").Value = -1362.683
$ws.Range("N132").Value = -13663.8329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 16113339
$ws.Range("I68").Value = 26026754
$ws.Range("K68").Value = 26026754
$ws.Range("M68").Value = -26026005
$ws.Range("H71").Value = 16113339
$ws.Range("I71").Value = 26026754
$ws.Range("K71").Value = 130133770
$ws.Range("M71").Value = -130130026
$ws.Range("H82").Value = 9093826
$ws.Range("I82").Value = 18183150
$ws.Range("J82").Value = 4500.6
$ws.Range("K82").Value = 18183150
$ws.Range("L82").Value = 4500.6
$ws.Range("M82").Value = -18182789
$ws.Range("N82").Value = -5222.6
$ws.Range("H85").Value = 9093826
$ws.Range("I85").Value = 18183150
$ws.Range("J85").Value = 4500.6
$ws.Range("K85").Value = 18183150
$ws.Range("L85").Value = 4500.6
$ws.Range("M85").Value = -18181902
$ws.Range("N85").Value = -6996.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1908.8864
$ws.Range("I132").Value = 730.9429
$ws.Range("K132").Value = 2192.8287
$ws.Range("M132").Value = 337.1713
$ws.Range("H136").Value = 3061.549
$ws.Range("I136").Value = 907.5897
$ws.Range("J136").Value = 10061.917
$ws.Range("K136").Value = 2722.7691
$ws.Range("L136").Value = 30185.751
$ws.Range("M136").Value = -172.7691
$ws.Range("N136").Value = -35285.751
